# Fixed 2016 powerplant generation data
# - Column O holds gwh_2016 values; many rows had a placeholder 0 that is
#   corrected to the real (re-extracted) generation figure. Two rows
#   (30, 34) turn out to be unresolvable / not-available, matching the
#   existing #N/A already present in columns P/Q for those plants.
# - The P/Q header labels are shortened from "generation_gwh_2017/2018"
#   to "gwh_2017"/"gwh_2018" to match the other gwh_* column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (P1, Q1) ---------------------------------------------
$ws.Range("P1").Value = "gwh_2017"
$ws.Range("Q1").Value = "gwh_2018"

# --- Column O (gwh_2016) corrected generation values --------------------
$ws.Range("O24").Value = 195.88466099999999
$ws.Range("O25").Value = 335.92700000000002
$ws.Range("O26").Value = 480.96699999999998
$ws.Range("O27").Value = 308.96699999999998
$ws.Range("O28").Value = 906.99599999999998
$ws.Range("O29").Value = 374.892
$ws.Range("O30").Value = "#N/A"
$ws.Range("O31").Value = 10.46898
$ws.Range("O32").Value = 11.308044000000001
$ws.Range("O34").Value = "#N/A"
$ws.Range("O35").Value = 8564.3790000000008
$ws.Range("O36").Value = 3365.7460000000001
$ws.Range("O37").Value = 3055.3049999999998
$ws.Range("O38").Value = 1002.6660000000001
$ws.Range("O39").Value = 9987.7520000000004
$ws.Range("O40").Value = 1128.2940000000001
$ws.Range("O41").Value = 5089.6549999999997
$ws.Range("O42").Value = 21.730399999999999
$ws.Range("O43").Value = 5806.8509999999997
$ws.Range("O44").Value = 5239.7030000000004
$ws.Range("O45").Value = 66.207300000000004
$ws.Range("O46").Value = 93.440449999999998
$ws.Range("O47").Value = 42.496450000000003
$ws.Range("O49").Value = 226.96944999999999
$ws.Range("O50").Value = 61.421349999999997
$ws.Range("O51").Value = 287.943049999999
$ws.Range("O52").Value = 811.53195000000005
$ws.Range("O53").Value = 12.676299999999999
$ws.Range("O54").Value = 191.58725000000001
$ws.Range("O55").Value = 91.80865
$ws.Range("O56").Value = 124.8526
$ws.Range("O57").Value = 2916.67
$ws.Range("O64").Value = 62.973549999999904
$ws.Range("O65").Value = 23.830249999999999
$ws.Range("O66").Value = 19.730849999999901
$ws.Range("O67").Value = 20.48705
$ws.Range("O68").Value = 17.382649999999899
$ws.Range("O69").Value = 168.58285000000001
$ws.Range("O74").Value = 3256.8510000000001
$ws.Range("O75").Value = 4652.2794999999996
$ws.Range("O76").Value = 3212.7
$ws.Range("O77").Value = 5912.3959999999997
$ws.Range("O78").Value = 3482.79

# --- Cosmetic view-state tweaks left by the author's Excel session ------
# Column O widened slightly and the sheet scrolled/selected further down.
$ws.Columns.Item(15).ColumnWidth = 23
$excel.Goto($ws.Range("H16"), $true)
$ws.Range("P5").Select()
